$wb = $excel.ActiveWorkbook

# --- Sheet "SPN" (sheet1) ---
$ws1 = $wb.Worksheets.Item("SPN")

# Update J7 from "Pendente" to "Resolvido"
$ws1.Range("J7").Value = "Resolvido"

# Append new rows 10-13
$spnRows = @(
    @("SPN", "Higor Cruz",     2025, 3, "20/01/2025", "24/01/2025", 320403, "01/2025", "20/01/2025", "Resolvido", "Willian Jones"),
    @("SPN", "Higor Cruz",     2025, 3, "20/01/2025", "24/01/2025", 320607, "01/2025", "20/01/2025", "Pendente",  "Willian Jones"),
    @("SPN", "Luan Pierry",    2025, 3, "20/01/2025", "24/01/2025", 320215, "01/2025", "20/01/2025", "Resolvido", "Willian Jones"),
    @("SPN", "Nadir Baseggio", 2025, 3, "20/01/2025", "24/01/2025", 319786, "01/2025", "20/01/2025", "Pendente",  "Willian Jones")
)

$r = 10
foreach ($row in $spnRows) {
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws1.Cells.Item($r, $c + 1).Value = $row[$c]
    }
    $r++
}

# --- Sheet "ITI" (sheet2) ---
$ws2 = $wb.Worksheets.Item("ITI")

$itiRows = @(
    @("ITI", "Alana Neris",        2025, 3, "20/01/2025", "24/01/2025", 320974, "01/2025", "20/01/2025", "Resolvido", "Emerson Simette"),
    @("ITI", "Simette",            2025, 3, "20/01/2025", "24/01/2025", 320245, "01/2025", "20/01/2025", "Resolvido", "Emerson Simette"),
    @("ITI", "Erick da Silva",     2025, 3, "20/01/2025", "24/01/2025", 320780, "01/2025", "20/01/2025", "Resolvido", "Emerson Simette"),
    @("ITI", "Erick da Silva",     2025, 3, "20/01/2025", "24/01/2025", 320854, "01/2025", "20/01/2025", "Resolvido", "Emerson Simette"),
    @("ITI", "Erick da Silva",     2025, 3, "20/01/2025", "24/01/2025", 313182, "01/2025", "20/01/2025", "Resolvido", "Emerson Simette"),
    @("ITI", "Erick da Silva",     2025, 3, "20/01/2025", "24/01/2025", 320964, "01/2025", "20/01/2025", "Resolvido", "Emerson Simette"),
    @("ITI", "Erick da Silva",     2025, 3, "20/01/2025", "24/01/2025", 320990, "01/2025", "20/01/2025", "Resolvido", "Emerson Simette"),
    @("ITI", "Erick da Silva",     2025, 3, "20/01/2025", "24/01/2025", 321158, "01/2025", "20/01/2025", "Resolvido", "Emerson Simette"),
    @("ITI", "Felipe Nascimento",  2025, 3, "20/01/2025", "24/01/2025", 320312, "01/2025", "20/01/2025", "Pendente",  "Emerson Simette"),
    @("ITI", "Jorgenaldo Reis",    2025, 3, "20/01/2025", "24/01/2025", 320458, "01/2025", "20/01/2025", "Resolvido", "Emerson Simette"),
    @("ITI", "Jorgenaldo Reis",    2025, 3, "20/01/2025", "24/01/2025", 320546, "01/2025", "20/01/2025", "Resolvido", "Emerson Simette"),
    @("ITI", "Jorgenaldo Reis",    2025, 3, "20/01/2025", "24/01/2025", 320776, "01/2025", "20/01/2025", "Resolvido", "Emerson Simette"),
    @("ITI", "Jorgenaldo Reis",    2025, 3, "20/01/2025", "24/01/2025", 320861, "01/2025", "20/01/2025", "Resolvido", "Emerson Simette"),
    @("ITI", "Jorgenaldo Reis",    2025, 3, "20/01/2025", "24/01/2025", 320895, "01/2025", "20/01/2025", "Resolvido", "Emerson Simette"),
    @("ITI", "Jorgenaldo Reis",    2025, 3, "20/01/2025", "24/01/2025", 320916, "01/2025", "20/01/2025", "Resolvido", "Emerson Simette"),
    @("ITI", "Jorgenaldo Reis",    2025, 3, "20/01/2025", "24/01/2025", 320941, "01/2025", "20/01/2025", "Resolvido", "Emerson Simette"),
    @("ITI", "Jorgenaldo Reis",    2025, 3, "20/01/2025", "24/01/2025", 320942, "01/2025", "20/01/2025", "Resolvido", "Emerson Simette"),
    @("ITI", "Jorgenaldo Reis",    2025, 3, "20/01/2025", "24/01/2025", 320991, "01/2025", "20/01/2025", "Resolvido", "Emerson Simette"),
    @("ITI", "Jorgenaldo Reis",    2025, 3, "20/01/2025", "24/01/2025", 321068, "01/2025", "20/01/2025", "Resolvido", "Emerson Simette"),
    @("ITI", "Michel Pessoa",      2025, 3, "20/01/2025", "24/01/2025", 320583, "01/2025", "20/01/2025", "Resolvido", "Emerson Simette")
)

$r = 14
foreach ($row in $itiRows) {
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws2.Cells.Item($r, $c + 1).Value = $row[$c]
    }
    $r++
}
